$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "wildcard"
$ws.Range("E2").Value = "Exact Page Redirect"

$ws.Range("D3").Value = "partial"
$ws.Range("E3").Value = "Section/Folder Redirect"
